# Daily attendance processing - 2025-12-22 01:36:26
# Reorders the "Recorded By" (column G) values so that "System" is listed
# first in the comma-separated list of recorders, instead of last -
# except for rows recorded solely/first by "admin@admin.com", which are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }

    if ($val.EndsWith(", System") -and -not $val.StartsWith("admin@admin.com")) {
        $parts = $val -split ", "
        $others = $parts[0..($parts.Length - 2)]
        $newVal = "System, " + ($others -join ", ")
        $cell.Value = $newVal
    }
}
